$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value. Values are written as text (matching
# the source data, which stores prices/percentages as inline strings rather
# than numbers) by forcing a text number-format before assignment.
$updates = [ordered]@{
    "D2" = "301.22"
    "D3" = "35.15"
    "E3" = "-0.42%"
    "D4" = "5.044"
    "E4" = "-1.77%"
    "D5" = "0.07981"
    "E5" = "-1.81%"
    "D6" = "1.904"
    "E6" = "-10.49%"
    "D7" = "7.800"
    "E7" = "-2.02%"
    "E8" = "-2.40%"
    "E9" = "0.81%"
    "D10" = "0.9219"
    "E10" = "-0.78%"
    "D11" = "0.1277"
    "E11" = "26.12%"
    "D12" = "0.1843"
    "E12" = "-1.36%"
    "D13" = "0.09866"
    "E13" = "8.35%"
    "D14" = "0.03576"
    "E14" = "-0.89%"
    "E15" = "-0.57%"
    "D16" = "0.001386"
    "E16" = "-3.95%"
    "D17" = "0.005783"
    "E17" = "1.66%"
    "D18" = "3.504"
    "E18" = "1.11%"
    "E19" = "-0.35%"
    "D20" = "0.1302"
    "E20" = "-2.87%"
    "D21" = "5.054"
    "E21" = "-1.01%"
    "D22" = "0.2402"
    "E22" = "8.17%"
    "D23" = "0.04503"
    "E23" = "-1.20%"
    "D24" = "0.001216"
    "E24" = "-2.46%"
    "D25" = "0.004781"
    "E25" = "1.72%"
    "E26" = "0.02%"
    "D27" = "0.0003004"
    "D39" = "0.01868"
    "E39" = "-4.70%"
    "D40" = "0.04694"
    "E40" = "-3.23%"
    "D41" = "0.007501"
    "E41" = "-2.90%"
    "D42" = "0.01029"
    "E42" = "31.20%"
    "D43" = "0.1322"
    "E43" = "-5.04%"
    "E44" = "-1.84%"
    "E45" = "-9.62%"
    "D46" = "0.00006248"
    "E46" = "-6.03%"
    "E47" = "0.01%"
    "E48" = "63.44%"
    "E49" = "-12.48%"
    "D50" = "0.00002103"
    "E50" = "0.01%"
    "E51" = "0.01%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
